{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Data description\" Heading1 paragraph that starts the section\n// being removed (that heading plus the paragraph right after it, which\n// begins \"For the analysis, a selective list of 57 stations ...\").\nlet targetIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"Data description\") {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex !== -1) {\n  items[targetIndex].delete();\n  items[targetIndex + 1].delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Data description\" Heading1 paragraph that starts the section\n# being removed (that heading plus the paragraph right after it).\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Trim() -eq \"Data description\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -ge 1) {\n    # Remove the \"Data description\" heading paragraph.\n    $d.Paragraphs.Item($targetIndex).Range.Delete()\n\n    # The paragraph that follows (\"For the analysis, a selective list of 57\n    # stations ... posed earlier.\") now sits at the same index; remove it too.\n    $d.Paragraphs.Item($targetIndex).Range.Delete()\n}\n"}
